$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at the top: a title row, a subtitle row, and a
# thin spacer row - pushing the existing questionnaire rows down by 3.
$ws.Rows("1:3").Insert()

# Row 1: "Heading 1" style title bound to [[name]]
$ws.Range("A1").Value = "[[name]]"
$ws.Range("A1").Style = "Heading 1"
$ws.Rows(1).RowHeight = 25.5

# Row 2: "Heading 2" style subtitle bound to [[description]]
$ws.Range("A2").Value = "[[description]]"
$ws.Range("A2").Style = "Heading 2"
$ws.Rows(2).RowHeight = 18

# Row 3: thin spacer row left under the heading's thick bottom border
$ws.Rows(3).RowHeight = 15.75

$ws.Range("A2").Select() | Out-Null
